# Add "Week 16" innings data for summer 2022 (new Week 33 column -> AH)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the added week column
$ws.Range("AH1").Value = "Week 33"

# New innings-count data points for the players who played this week
$ws.Range("AH3").Value = 3.5
$ws.Range("AH6").Value = 10
$ws.Range("AH7").Value = 9
$ws.Range("AH8").Value = 5
$ws.Range("AH9").Value = 1.5
